$wb = $excel.ActiveWorkbook

# Sheet ALC, row 2 (Leve Item ID 5489)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value2 = 93.333336
$ws.Range("I2").Value2 = 93.333336
$ws.Range("J2").Value2 = 0
$ws.Range("K2").Value2 = 93.333336
$ws.Range("L2").Value2 = 0
$ws.Range("M2").Value2 = 19.666664
$ws.Range("N2").ClearContents() | Out-Null

# Sheet ALC, row 5 (Leve Item ID 5503)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value2 = 198.75
$ws.Range("I5").Value2 = 100
$ws.Range("K5").Value2 = 100
$ws.Range("M5").Value2 = 15

# Sheet ALC, row 9 (Leve Item ID 5487)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value2 = 62.42857
$ws.Range("I9").Value2 = 62.42857
$ws.Range("K9").Value2 = 62.42857
$ws.Range("M9").Value2 = 106.57143

# Sheet ALC, row 18 (Leve Item ID 5471)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value2 = 287.25
$ws.Range("I18").Value2 = 287.25
$ws.Range("K18").Value2 = 287.25
$ws.Range("M18").Value2 = -3.25

# Sheet ALC, row 40 (Leve Item ID 5505)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value2 = 500
$ws.Range("I40").Value2 = 500
$ws.Range("K40").Value2 = 500
$ws.Range("M40").Value2 = -325

# Sheet ALC, row 43 (Leve Item ID 5472)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value2 = 4266.6665
$ws.Range("J43").Value2 = 1400
$ws.Range("L43").Value2 = 1400
$ws.Range("N43").Value2 = -1538

# Sheet ALC, row 70 (Leve Item ID 12604)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value2 = 4500
$ws.Range("I70").Value2 = 0
$ws.Range("J70").Value2 = 4500
$ws.Range("K70").Value2 = 0
$ws.Range("L70").Value2 = 13500
$ws.Range("M70").ClearContents() | Out-Null
$ws.Range("N70").Value2 = -14040

# Sheet ALC, row 73 (Leve Item ID 12604)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value2 = 4500
$ws.Range("I73").Value2 = 0
$ws.Range("J73").Value2 = 4500
$ws.Range("K73").Value2 = 0
$ws.Range("L73").Value2 = 13500
$ws.Range("M73").ClearContents() | Out-Null
$ws.Range("N73").Value2 = -15372

# Sheet ALC, row 80 (Leve Item ID 12605)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value2 = 838
$ws.Range("I80").Value2 = 250
$ws.Range("J80").Value2 = 985
$ws.Range("K80").Value2 = 750
$ws.Range("L80").Value2 = 2955
$ws.Range("M80").Value2 = 248
$ws.Range("N80").Value2 = -4951

# Sheet ALC, row 83 (Leve Item ID 12605)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value2 = 838
$ws.Range("I83").Value2 = 250
$ws.Range("J83").Value2 = 985
$ws.Range("K83").Value2 = 2250
$ws.Range("L83").Value2 = 8865
$ws.Range("M83").Value2 = 2742
$ws.Range("N83").Value2 = -18849

# Sheet ALC, row 138 (Leve Item ID 44169)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value2 = 875
$ws.Range("I138").Value2 = 875
$ws.Range("K138").Value2 = 2625
$ws.Range("M138").Value2 = 2515

# Sheet ARM, row 4 (Leve Item ID 5071)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value2 = 156.36363
$ws.Range("I4").Value2 = 162.1
$ws.Range("K4").Value2 = 162.1
$ws.Range("M4").Value2 = -46.09999999999999

# Sheet ARM, row 5 (Leve Item ID 5091)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value2 = 175
$ws.Range("I5").Value2 = 175
$ws.Range("K5").Value2 = 175
$ws.Range("M5").Value2 = -63

# Sheet ARM, row 38 (Leve Item ID 2260)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H38").Value2 = 0
$ws.Range("I38").Value2 = 0
$ws.Range("K38").Value2 = 0
$ws.Range("M38").ClearContents() | Out-Null

# Sheet ARM, row 63 (Leve Item ID 12528)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value2 = 7589.6
$ws.Range("I63").Value2 = 8483.333000000001
$ws.Range("K63").Value2 = 8483.333000000001
$ws.Range("M63").Value2 = -7797.333000000001

# Sheet ARM, row 66 (Leve Item ID 12528)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value2 = 7589.6
$ws.Range("I66").Value2 = 8483.333000000001
$ws.Range("K66").Value2 = 42416.665
$ws.Range("M66").Value2 = -38984.665

# Sheet ARM, row 88 (Leve Item ID 12530)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value2 = 800
$ws.Range("I88").Value2 = 1000
$ws.Range("K88").Value2 = 1000
$ws.Range("M88").Value2 = -594

# Sheet ARM, row 91 (Leve Item ID 12530)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value2 = 800
$ws.Range("I91").Value2 = 1000
$ws.Range("K91").Value2 = 1000
$ws.Range("M91").Value2 = 404

# Sheet ARM, row 101 (Leve Item ID 18518)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H101").Value2 = 37777
$ws.Range("J101").Value2 = 37777
$ws.Range("L101").Value2 = 37777
$ws.Range("N101").Value2 = -44267

# Sheet BSM, row 4 (Leve Item ID 5091)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value2 = 175
$ws.Range("I4").Value2 = 175
$ws.Range("K4").Value2 = 175
$ws.Range("M4").Value2 = -60

# Sheet BSM, row 86 (Leve Item ID 12526)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value2 = 1050
$ws.Range("I86").Value2 = 900
$ws.Range("J86").Value2 = 1500
$ws.Range("K86").Value2 = 900
$ws.Range("L86").Value2 = 1500
$ws.Range("M86").Value2 = 223
$ws.Range("N86").Value2 = -3746

# Sheet BSM, row 89 (Leve Item ID 12526)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value2 = 1050
$ws.Range("I89").Value2 = 900
$ws.Range("J89").Value2 = 1500
$ws.Range("K89").Value2 = 4500
$ws.Range("L89").Value2 = 7500
$ws.Range("M89").Value2 = 1116
$ws.Range("N89").Value2 = -18732

# Sheet CRP, row 7 (Leve Item ID 5361)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value2 = 45.57143
$ws.Range("I7").Value2 = 33.8
$ws.Range("J7").Value2 = 75
$ws.Range("K7").Value2 = 33.8
$ws.Range("L7").Value2 = 75
$ws.Range("M7").Value2 = 79.2
$ws.Range("N7").Value2 = -301

# Sheet CRP, row 22 (Leve Item ID 5367)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value2 = 339.14285
$ws.Range("I22").Value2 = 314.8
$ws.Range("K22").Value2 = 314.8
$ws.Range("M22").Value2 = 35.19999999999999

# Sheet CRP, row 62 (Leve Item ID 12580)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value2 = 3500
$ws.Range("I62").Value2 = 0
$ws.Range("J62").Value2 = 3500
$ws.Range("K62").Value2 = 0
$ws.Range("L62").Value2 = 3500
$ws.Range("M62").ClearContents() | Out-Null
$ws.Range("N62").Value2 = -4748

# Sheet CRP, row 65 (Leve Item ID 12580)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value2 = 3500
$ws.Range("I65").Value2 = 0
$ws.Range("J65").Value2 = 3500
$ws.Range("K65").Value2 = 0
$ws.Range("L65").Value2 = 17500
$ws.Range("M65").ClearContents() | Out-Null
$ws.Range("N65").Value2 = -23740

# Sheet CRP, row 88 (Leve Item ID 10608)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H88").Value2 = 14850.929
$ws.Range("J88").Value2 = 14850.929
$ws.Range("L88").Value2 = 14850.929
$ws.Range("N88").Value2 = -15662.929

# Sheet CRP, row 91 (Leve Item ID 10608)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H91").Value2 = 14850.929
$ws.Range("J91").Value2 = 14850.929
$ws.Range("L91").Value2 = 14850.929
$ws.Range("N91").Value2 = -17658.929

# Sheet CRP, row 134 (Leve Item ID 44020)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value2 = 1997.5
$ws.Range("J134").Value2 = 1995
$ws.Range("L134").Value2 = 5985
$ws.Range("N134").Value2 = -11055

# Sheet CUL, row 12 (Leve Item ID 4854)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value2 = 171.83333
$ws.Range("I12").Value2 = 48.4
$ws.Range("K12").Value2 = 145.2
$ws.Range("M12").Value2 = 27.80000000000001

# Sheet CUL, row 23 (Leve Item ID 4858)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value2 = 381
$ws.Range("I23").Value2 = 50
$ws.Range("J23").Value2 = 491.33334
$ws.Range("K23").Value2 = 150
$ws.Range("L23").Value2 = 1474.00002
$ws.Range("M23").Value2 = 85
$ws.Range("N23").Value2 = -1944.00002

# Sheet CUL, row 108 (Leve Item ID 27853)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H108").Value2 = 578.3333
$ws.Range("J108").Value2 = 1335
$ws.Range("L108").Value2 = 4005
$ws.Range("N108").Value2 = -9765

# Sheet GSM, row 2 (Leve Item ID 5062)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value2 = 80.333336
$ws.Range("I2").Value2 = 80.333336
$ws.Range("K2").Value2 = 80.333336
$ws.Range("M2").Value2 = 32.666664

# Sheet GSM, row 43 (Leve Item ID 4218)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value2 = 7985
$ws.Range("I43").Value2 = 5001
$ws.Range("J43").Value2 = 10969
$ws.Range("K43").Value2 = 5001
$ws.Range("L43").Value2 = 10969
$ws.Range("M43").Value2 = -4850
$ws.Range("N43").Value2 = -11271

# Sheet GSM, row 80 (Leve Item ID 12521)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value2 = 2493.5
$ws.Range("I80").Value2 = 2493.5
$ws.Range("J80").Value2 = 0
$ws.Range("K80").Value2 = 2493.5
$ws.Range("L80").Value2 = 0
$ws.Range("M80").Value2 = -1495.5
$ws.Range("N80").ClearContents() | Out-Null

# Sheet GSM, row 83 (Leve Item ID 12521)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value2 = 2493.5
$ws.Range("I83").Value2 = 2493.5
$ws.Range("J83").Value2 = 0
$ws.Range("K83").Value2 = 12467.5
$ws.Range("L83").Value2 = 0
$ws.Range("M83").Value2 = -7475.5
$ws.Range("N83").ClearContents() | Out-Null

# Sheet GSM, row 123 (Leve Item ID 34150)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value2 = 46265.875
$ws.Range("J123").Value2 = 46265.875
$ws.Range("L123").Value2 = 46265.875
$ws.Range("N123").Value2 = -51165.875

# Sheet LTW, row 22 (Leve Item ID 5277)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value2 = 5399.5
$ws.Range("J22").Value2 = 5399.5
$ws.Range("L22").Value2 = 5399.5
$ws.Range("N22").Value2 = -5989.5

# Sheet LTW, row 27 (Leve Item ID 5277)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value2 = 5399.5
$ws.Range("J27").Value2 = 5399.5
$ws.Range("L27").Value2 = 5399.5
$ws.Range("N27").Value2 = -5613.5

# Sheet LTW, row 46 (Leve Item ID 5282)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value2 = 0
$ws.Range("I46").Value2 = 0
$ws.Range("J46").Value2 = 0
$ws.Range("K46").Value2 = 0
$ws.Range("L46").Value2 = 0
$ws.Range("M46").ClearContents() | Out-Null
$ws.Range("N46").ClearContents() | Out-Null

# Sheet LTW, row 82 (Leve Item ID 12565)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value2 = 1666.3334
$ws.Range("I82").Value2 = 1000
$ws.Range("K82").Value2 = 1000
$ws.Range("M82").Value2 = -639

# Sheet LTW, row 85 (Leve Item ID 12565)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value2 = 1666.3334
$ws.Range("I85").Value2 = 1000
$ws.Range("K85").Value2 = 1000
$ws.Range("M85").Value2 = 248

# Sheet LTW, row 100 (Leve Item ID 19995)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value2 = 1850
$ws.Range("I100").Value2 = 1850
$ws.Range("K100").Value2 = 1850
$ws.Range("M100").Value2 = -1309

# Sheet LTW, row 122 (Leve Item ID 36247)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value2 = 3766.5557
$ws.Range("I122").Value2 = 2999.5
$ws.Range("K122").Value2 = 8998.5
$ws.Range("M122").Value2 = -6548.5

# Sheet WVR, row 4 (Leve Item ID 2996)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value2 = 2067
$ws.Range("I4").Value2 = 602
$ws.Range("J4").Value2 = 2799.5
$ws.Range("K4").Value2 = 602
$ws.Range("L4").Value2 = 2799.5
$ws.Range("M4").Value2 = -489
$ws.Range("N4").Value2 = -3025.5

# Sheet WVR, row 62 (Leve Item ID 12589)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value2 = 2875
$ws.Range("I62").Value2 = 2500
$ws.Range("J62").Value2 = 3250
$ws.Range("K62").Value2 = 2500
$ws.Range("L62").Value2 = 3250
$ws.Range("M62").Value2 = -1876
$ws.Range("N62").Value2 = -4498

# Sheet WVR, row 65 (Leve Item ID 12589)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value2 = 2875
$ws.Range("I65").Value2 = 2500
$ws.Range("J65").Value2 = 3250
$ws.Range("K65").Value2 = 12500
$ws.Range("L65").Value2 = 16250
$ws.Range("M65").Value2 = -9380
$ws.Range("N65").Value2 = -22490

# Sheet WVR, row 81 (Leve Item ID 12596)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value2 = 400
$ws.Range("I81").Value2 = 400
$ws.Range("J81").Value2 = 0
$ws.Range("K81").Value2 = 800
$ws.Range("L81").Value2 = 0
$ws.Range("M81").Value2 = 261
$ws.Range("N81").ClearContents() | Out-Null

# Sheet WVR, row 84 (Leve Item ID 12596)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value2 = 400
$ws.Range("I84").Value2 = 400
$ws.Range("J84").Value2 = 0
$ws.Range("K84").Value2 = 4000
$ws.Range("L84").Value2 = 0
$ws.Range("M84").Value2 = 1304
$ws.Range("N84").ClearContents() | Out-Null

# Sheet WVR, row 132 (Leve Item ID 44029)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value2 = 1784.8572
$ws.Range("I132").Value2 = 1784.8572
$ws.Range("K132").Value2 = 5354.571599999999
$ws.Range("M132").Value2 = -2824.571599999999
